$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A36").Value = "Emanuele Miorandi"
$ws.Range("B36").Value = "Alberto Cerisara | SHARK ATTACK"
$ws.Range("C36").Value = "Alessio Bragagna | SHARK ATTACK"
$ws.Range("D36").Value = "Filippo Benetti | I Magnifici"
$ws.Range("E36").Value = "Alessandro Maffei | FC SAVIGNANO"
$ws.Range("F36").Value = "Moris Benedetti | Gli Introvabili"
